$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 150-151, shifting existing rows 150-178 down to 152-180
$ws.Rows("150:151").Insert()

# Row 150
$ws.Range("A150").Value = 4
$ws.Range("B150").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C150").Value = 'Los Lagos'
$ws.Range("D150").Value = 44642
$ws.Range("E150").Value = 10
$ws.Range("F150").Value = 'Fruta'
$ws.Range("G150").Value = 100103
$ws.Range("H150").Value = 'Frutos de hueso (carozo)'
$ws.Range("I150").Value = 100103002
$ws.Range("J150").Value = 'Ciruela'
$ws.Range("K150").Value = 'Angeleno'
$ws.Range("L150").Value = 'Primera'
$ws.Range("M150").Value = 400
$ws.Range("N150").Value = 14000
$ws.Range("O150").Value = 15000
$ws.Range("P150").Value = 14500
$ws.Range("Q150").Value = '$/caja 15 kilos granel'
$ws.Range("R150").Value = 'Provincia de Curicó'
$ws.Range("S150").Value = 967
$ws.Range("T150").Value = 15

# Row 151
$ws.Range("A151").Value = 4
$ws.Range("B151").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C151").Value = 'Los Lagos'
$ws.Range("D151").Value = 44642
$ws.Range("E151").Value = 10
$ws.Range("F151").Value = 'Fruta'
$ws.Range("G151").Value = 100103
$ws.Range("H151").Value = 'Frutos de hueso (carozo)'
$ws.Range("I151").Value = 100103002
$ws.Range("J151").Value = 'Ciruela'
$ws.Range("K151").Value = 'Angeleno'
$ws.Range("L151").Value = 'Segunda'
$ws.Range("M151").Value = 200
$ws.Range("N151").Value = 13000
$ws.Range("O151").Value = 13000
$ws.Range("P151").Value = 13000
$ws.Range("Q151").Value = '$/caja 15 kilos granel'
$ws.Range("R151").Value = 'Provincia de Curicó'
$ws.Range("S151").Value = 867
$ws.Range("T151").Value = 15

